$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Directeur" column (D) with manager names per branch row.
# Order below mirrors the original authoring order (shared-string insertion
# order), so the rebuilt sharedStrings table lines up cell-for-cell with
# the source workbook.
$ws.Range("D2").Value = "Martin Hudon"
$ws.Range("D3").Value = "Eric Vanier"
$ws.Range("D4").Value = "Caroline Cyr"
$ws.Range("D5").Value = "Yannick Lazare"
$ws.Range("D6").Value = "Simon Huard"
$ws.Range("D7").Value = "Bruno Adam"
$ws.Range("D8").Value = "Claudio Fazioli"
$ws.Range("D9").Value = "Eric Vanier"
$ws.Range("D10").Value = "Yannick Lazare"
$ws.Range("D11").Value = "Simon Huard"
$ws.Range("D12").Value = "Claudio Fazioli"
$ws.Range("D13").Value = "Karolane Roy"
$ws.Range("D14").Value = "Karolane Roy"
$ws.Range("D15").Value = "Caroline Cyr"
$ws.Range("D20").Value = "Ludovic Gérard"
$ws.Range("D22").Value = "Chantal Maltais"
$ws.Range("D21").Value = "Danny Pronovost"
$ws.Range("D23").Value = "Yannick Blanchet"
$ws.Range("D24").Value = "Yannick Blanchet"
$ws.Range("D25").Value = "Eric Savard"
$ws.Range("D26").Value = "Danny Pronovost"
$ws.Range("D1").Value = "Directeur"

# New column width for the added Directeur column (stored width ends up
# as 16 once Excel applies its fixed character-to-width padding).
$ws.Columns.Item(4).ColumnWidth = 15.1666666666667

# Selection ends up on D11 in the saved file.
$ws.Range("D11").Select()
